# Flip the "Execute" column to "Y" for all test-case rows on the
# "Scenarios" sheet (rows 2-13), enabling the critical-path test suite
# (including the newly added login test) to run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

$ws.Range("A2:A13").Value = "Y"
